$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.847.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.48%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.829.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.61%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9983"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.28%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'244.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.6908"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.99%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.9994"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.16%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.07665"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.22%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.3046"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.36%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -3.31%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07795"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.10%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'93.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.92%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.829.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.03%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.086"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.90%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.6803"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.62%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'6.441"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.61%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.000008238"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.77%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'28.816.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.33%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'241.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.63%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'2.073.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.48%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'12.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.06%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.06%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'7.449"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.06%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.15%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.1497"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.29%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'161.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.11%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'8.735"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.98%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'18.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.21%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.541"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.21%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.229"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.79%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.165"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.08%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.190"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.16%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.05118"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.02%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.7744"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +3.30%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.856"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.00%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.139"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.98%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.14%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'Maker"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'1.269.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.83%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'VeChain"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'0.01852"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.63%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.697"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.70%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.9577"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +6.38%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'6.068"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.28%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'106.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.84%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.9998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.13%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'9.652"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.19%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.5162"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.40%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.972.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.69%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'64.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -7.13%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'RenderToken"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'1.750"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.46%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'BabyDogeCoin"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.00000000119"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -4.42%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'6.940"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.67%  "
$ws.Range("E51").Style = "Normal"
